# "clean-up of input tables"
# The only substantive, COM-reachable change captured by the diff is the
# worksheet being renamed from "updated" to "Tabelle1".
# (The remaining hunks in the diff -- the x15ac:absPath save location,
# the bookViews window geometry, sheetFormatPr defaults, the per-row
# x14ac:dyDescent hints and the column B/C width rounding -- are simply
# re-save fingerprints left by the Excel version/machine that produced
# the file and are not data the object model exposes for editing, so
# they are intentionally left untouched here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Tabelle1"
